$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": append a new date column AB ("11-jul") with its values.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Write the new header value first, then clone the header formatting from the
# existing last column (AA1) via a formats-only paste so the new cell reuses
# the same cell style as the rest of the header row.
$ws1.Range("AB1").Value = "11-jul"
$ws1.Range("AA1").Copy()
$ws1.Range("AB1").PasteSpecial(-4122)

$ws1.Range("AB2").Value = 73.12
$ws1.Range("AB3").Value = 62.48
$ws1.Range("AB4").Value = 51.1
$ws1.Range("AB5").Value = 50.61
$ws1.Range("AB6").Value = 44.16
$ws1.Range("AB7").Value = 39.78
$ws1.Range("AB8").Value = 56.01
$ws1.Range("AB9").Value = 65.81999999999999
$ws1.Range("AB10").Value = 57.2
$ws1.Range("AB11").Value = 60.05
$ws1.Range("AB12").Value = 50
$ws1.Range("AB13").Value = 28.35
$ws1.Range("AB14").Value = 39.46
$ws1.Range("AB15").Value = 28.11
$ws1.Range("AB16").Value = 32.04
$ws1.Range("AB17").Value = 22.38
$ws1.Range("AB18").Value = 30.7
$ws1.Range("AB19").Value = 47.85
$ws1.Range("AB20").Value = 66.56999999999999
$ws1.Range("AB21").Value = 77.7
$ws1.Range("AB22").Value = 71.92
$ws1.Range("AB23").Value = 61.92
$ws1.Range("AB24").Value = 98.93000000000001
$ws1.Range("AB25").Value = 89.98999999999999

# ---------------------------------------------------------------------------
# Sheet "Gaz": append a new row 25 for 2025-07-09.
#
# The date-like text must stay literal text (matching every other "Date"
# cell in column A), not get auto-converted to a date serial number. Force
# text entry via a temporary Text number format, then drop that formatting
# again so the cell ends up with no explicit style, just like its neighbours.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A25").NumberFormat = "@"
$ws2.Range("A25").Value = "2025-07-09"
$ws2.Range("A25").ClearFormats()
$ws2.Range("B25").Value = 33.6

# ---------------------------------------------------------------------------
# Sheet "CO2": append a new row 25 for 2025-07-09.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A25").NumberFormat = "@"
$ws3.Range("A25").Value = "2025-07-09"
$ws3.Range("A25").ClearFormats()
$ws3.Range("B25").Value = 69.65000000000001
